# Backlog update: Prasso now saves results.
#
# The "As a user / I want my results saved" backlog item (originally the
# top-priority row, row 2) is now done: Prasso saves results. So its
# Story Points drop to 0 and its Status becomes "Done". The sheet is kept
# sorted descending by Story Points (column C), so that row slides down to
# just above the existing block of completed (0-point) items - i.e. every
# row between the old and new position shifts up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$firstRow = 2
$lastRow  = 23

# Snapshot the current A/B/C/E values for rows 2..23 before overwriting
# anything (columns: A=ID, B=Description, C=Story Points, E=Status).
$colsToShift = 1, 2, 3, 5
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    foreach ($c in $colsToShift) {
        $key = [string]$r + "_" + [string]$c
        $snapshot[$key] = $ws.Cells.Item($r, $c).Value2
    }
}

# Shift rows (firstRow+1)..lastRow up by one row.
for ($r = $firstRow; $r -lt $lastRow; $r++) {
    foreach ($c in $colsToShift) {
        $srcKey = [string]($r + 1) + "_" + [string]$c
        $ws.Cells.Item($r, $c).Value = $snapshot[$srcKey]
    }
}

# The row that fell off the top (old row 2 = "I want my results saved")
# is now finished, so it is placed at the bottom of the block (row 23)
# with 0 remaining story points and status "Done".
$ws.Cells.Item($lastRow, 1).Value = $snapshot[[string]$firstRow + "_1"]
$ws.Cells.Item($lastRow, 2).Value = $snapshot[[string]$firstRow + "_2"]
$ws.Cells.Item($lastRow, 3).Value = 0
$ws.Cells.Item($lastRow, 5).Value = "Done"

# Match the author's final cursor position (cell E2 was left selected).
[void]$ws.Range("E2").Select()
